$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:45:30"
$ws1.Range("A3").Value = "Total filas: 15"

$ws1.Range("A6").Value  = "04:45:30"
$ws1.Range("B6").Value  = "04:46"
$ws1.Range("D6").Value  = 1

$ws1.Range("A7").Value  = "04:45:30"
$ws1.Range("D7").Value  = 8

$ws1.Range("A8").Value  = "04:45:30"
$ws1.Range("D8").Value  = 26

$ws1.Range("A9").Value  = "04:45:30"
$ws1.Range("D9").Value  = 37

$ws1.Range("A10").Value = "04:45:30"
$ws1.Range("B10").Value = "05:31"
$ws1.Range("D10").Value = 46

$ws1.Range("A11").Value = "04:45:30"
$ws1.Range("D11").Value = 59

$ws1.Range("A12").Value = "04:45:30"
$ws1.Range("D12").Value = 67

$ws1.Range("A13").Value = "04:45:30"
$ws1.Range("D13").Value = 76

$ws1.Range("A14").Value = "04:45:30"
$ws1.Range("B14").Value = "06:03"
$ws1.Range("D14").Value = 78

$ws1.Range("A15").Value = "04:45:30"
$ws1.Range("D15").Value = 86

$ws1.Range("A16").Value = "04:45:30"
$ws1.Range("D16").Value = 99

$ws1.Range("A17").Value = "04:45:30"
$ws1.Range("D17").Value = 102

$ws1.Range("A18").Value = "04:45:30"
$ws1.Range("D18").Value = 106

$ws1.Range("A19").Value = "04:45:30"
$ws1.Range("D19").Value = 106

# New row 20
$ws1.Range("A20").Value = "04:45:30"
$ws1.Range("B20").Value = "06:39"
$ws1.Range("C20").Value = "225_C ROCA-H SUR"
$ws1.Range("D20").Value = 114
$ws1.Range("E20").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:45:30"

$ws2.Range("A6").Value = "04:45:30"
$ws2.Range("B6").Value = "04:46"
$ws2.Range("D6").Value = 1

$ws2.Range("A7").Value = "04:45:30"
$ws2.Range("D7").Value = 86

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:45:30"
